# Atualização automática via cronjob
# Remove the oldest day's rows (2025-05-21) from the "vendas atipicas" sheet
# and refresh the remaining rows' data, shifting everything up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows belonging to 2025-05-21 (rows 2 and 3); the rest
# of the table shifts up automatically.
$ws.Range("A2:A3").EntireRow.Delete()

# Final dataset (10 data rows) after the cronjob refresh.
$data = @(
    @(6, "2025-05-22", 33,  "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA",       "000897", "AGUA SANITARIA GLOBO SAN 5L",                                    29,   $false),
    @(7, "2025-05-22", 10,  "AMAZONIA REFEICOES E SERVICOS LTDA",                    "000999", "ESCOVA DE ACO 3 FILEIRAS",                                       0,    $false),
    @(8, "2025-05-22", 4,   "AMAZONIA REFEICOES E SERVICOS LTDA",                    "000426", "KIT DESCASCADOR DE LEGUMES KEITA",                               4,    $false),
    @(0, "2025-05-26", 200, "MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.","000098", "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM",              1791, $false),
    @(2, "2025-05-26", 40,  "MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.","000288", "TOUCA DESCARTAVEL TNT TALGE PCT C/ 100 UND",                     303,  $false),
    @(4, "2025-05-26", 96,  "MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.","000032", "LIMPADOR VEJA MULTIUSO GOLD 500ML",                              1062, $false),
    @(1, "2025-05-27", 50,  "AMAZONIA REFEICOES E SERVICOS LTDA",                    "000103", "AVENTAL PVC FORRADO PLUS 1,20X0,65 BRANCO C.A. 28303 BRASCAMP",  68,   $false),
    @(3, "2025-05-27", 24,  "AMAZONIA REFEICOES E SERVICOS LTDA",                    "001023", "FILME PVC 30X8X500MT",                                           -12,  $true),
    @(5, "2025-05-27", 40,  "AMAZONIA REFEICOES E SERVICOS LTDA",                    "000855", "SACO PLAST BD 50X80 N30 PC/25",                                  -8,   $false),
    @(9, "2025-05-27", 26,  "AMAZONIA REFEICOES E SERVICOS LTDA",                    "000152", "COPO DESCARTAVEL BRANCO CRISTALCOPO 180ML CX C\25",              -47,  $false)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]

    # Column B ("Dia") holds an ISO date string like "2025-05-22" — force
    # text storage so Excel doesn't auto-convert it to a date serial.
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $r[1]
    $bCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]

    # Column E ("id_produto") holds zero-padded codes like "000897" — force
    # text storage so Excel doesn't drop the leading zeros as a number.
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $r[4]
    $eCell.Style = "Normal"

    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}
